$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.177197694778442
$ws.Range("B1").Value = 2.420051574707031
$ws.Range("D1").Value = 2.332631587982178
$ws.Range("E1").Value = 1.202049732208252
